# Auto-generated edit script applying cryptos.xlsx price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.598.62"
$ws.Range("E2").Value = "  -1.46%  "
$ws.Range("D3").Value = "2.545.52"
$ws.Range("E3").Value = "  -0.71%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.13"
$ws.Range("E5").Value = "  -1.86%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "99.74"
$ws.Range("E6").Value = "  +2.43%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.569"
$ws.Range("E7").Value = "  -1.21%  "
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.527"
$ws.Range("E9").Value = "  -2.72%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.99"
$ws.Range("E10").Value = "  +0.42%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0803"
$ws.Range("E11").Value = "  -1.33%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.35"
$ws.Range("E12").Value = "  -2.34%  "
$ws.Range("E13").Value = "  -0.08%  "
$ws.Range("D14").Value = "2.944.24"
$ws.Range("E14").Value = "  -0.43%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "16.26"
$ws.Range("E15").Value = "  +7.46%  "
$ws.Range("D16").Value = "2.517.63"
$ws.Range("E16").Value = "  +1.53%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.839"
$ws.Range("E17").Value = "  -1.38%  "
$ws.Range("D18").Value = "42.595.37"
$ws.Range("E18").Value = "  -1.34%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.77"
$ws.Range("E19").Value = "  -1.17%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.31"
$ws.Range("E20").Value = "  -2.84%  "
$ws.Range("D21").Value = "0.0₃0951"
$ws.Range("E21").Value = "  -1.61%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "69.05"
$ws.Range("E22").Value = "  -1.44%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "242.27"
$ws.Range("E23").Value = "  -4.66%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.89"
$ws.Range("E24").Value = "  -2.54%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.06"
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.37"
$ws.Range("E26").Value = "  -2.41%  "
$ws.Range("E28").Value = "  -2.69%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "39.86"
$ws.Range("E29").Value = "  -2.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.08"
$ws.Range("E30").Value = "  -2.24%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "159.22"
$ws.Range("E31").Value = "  +1.92%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.67"
$ws.Range("E32").Value = "  -3.74%  "
$ws.Range("E33").Value = "  +13.85%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0798"
$ws.Range("E34").Value = "  -0.32%  "
$ws.Range("B35").Value = "WEMIXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.62"
$ws.Range("E35").Value = "  -3.23%  "
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.05"
$ws.Range("E36").Value = "  -3.99%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.17"
$ws.Range("E37").Value = "  -5.22%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "17.99"
$ws.Range("E38").Value = "  -6.82%  "
$ws.Range("E39").Value = "  -1.06%  "
$ws.Range("E40").Value = "  -0.64%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.22"
$ws.Range("E41").Value = "  +9.97%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "21.55"
$ws.Range("E42").Value = "  -3.40%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.35"
$ws.Range("E43").Value = "  +2.81%  "
$ws.Range("E44").Value = "  +0.26%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0296"
$ws.Range("E45").Value = "  -3.00%  "
$ws.Range("D46").Value = "1.959.33"
$ws.Range("E46").Value = "  -1.96%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.91"
$ws.Range("E47").Value = "  -2.16%  "
$ws.Range("D48").Value = "2.801.09"
$ws.Range("E48").Value = "  -0.25%  "
$ws.Range("B49").Value = "BitcoinSV"
$ws.Range("C49").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "80.74"
$ws.Range("E49").Value = "  -5.49%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.192"
$ws.Range("E50").Value = "  -0.46%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.39"
$ws.Range("E51").Value = "  -3.41%  "
